$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 0.8
$ws.Range("F6").Value = 0.8272727272727273
$ws.Range("G14").Value = 0.7097744360902256
$ws.Range("G15").Value = 0.6947368421052631
$ws.Range("G16").Value = 0.6947368421052631
$ws.Range("F17").Value = 0.5669172932330827
$ws.Range("G17").Value = 0.6962406015037594
$ws.Range("G18").Value = 0.9055555555555556
$ws.Range("F19").Value = 0.8044444444444444
$ws.Range("G21").Value = 0.8722499999999996
$ws.Range("F22").Value = 0.8934782608695653
$ws.Range("G23").Value = 0.9630434782608696
$ws.Range("F24").Value = 0.8716847826086956
$ws.Range("G24").Value = 0.9652173913043478
$ws.Range("F25").Value = 0.8695652173913043
$ws.Range("F42").Value = 0.5848484848484848
$ws.Range("G42").Value = 0.7075757575757575
$ws.Range("F43").Value = 0.5696969696969697
$ws.Range("G43").Value = 0.693939393939394
$ws.Range("F44").Value = 0.5666666666666667
$ws.Range("G44").Value = 0.6893939393939394
$ws.Range("G45").Value = 0.693939393939394
$ws.Range("F50").Value = 0.7850746268656716
$ws.Range("F51").Value = 0.7462686567164178
$ws.Range("F52").Value = 0.7343283582089553
$ws.Range("G53").Value = 0.8716417910447761
$ws.Range("G54").Value = 0.9066666666666666
$ws.Range("G55").Value = 0.8853333333333333
$ws.Range("G56").Value = 0.8826666666666667
$ws.Range("F57").Value = 0.7466666666666666
$ws.Range("G57").Value = 0.8746666666666667
$ws.Range("F62").Value = 0.8634920634920635
$ws.Range("F70").Value = 0.7948164146868251
$ws.Range("G70").Value = 0.8479481641468682
$ws.Range("F71").Value = 0.7697624190064795
$ws.Range("F72").Value = 0.7680345572354211
$ws.Range("G72").Value = 0.824622030237581
$ws.Range("G73").Value = 0.8215982721382289
